$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 79243
$ws.Range("B3").Value = 79243
$ws.Range("B4").Value = 91828
$ws.Range("B5").Value = 79243
$ws.Range("B6").Value = 79243
$ws.Range("B7").Value = 79243
$ws.Range("B8").Value = 79243
$ws.Range("B9").Value = 79243
$ws.Range("B10").Value = 79243
$ws.Range("B11").Value = 79243
$ws.Range("B12").Value = 79243
$ws.Range("B13").Value = 79243
$ws.Range("B14").Value = 79243
$ws.Range("B15").Value = 91828
$ws.Range("A16").Value = 130865706
$ws.Range("B16").Value = 79243
$ws.Range("Q16").Value = 447322
$ws.Range("R16").Value = 7042723
$ws.Range("S16").Value = 7
$ws.Range("Z16").Value = '14:33'
$ws.Range("AB16").Value = '14:33'
$ws.Range("A17").Value = 130865707
$ws.Range("B17").Value = 79243
$ws.Range("Q17").Value = 447282
$ws.Range("R17").Value = 7042740
$ws.Range("S17").Value = 9
$ws.Range("Z17").Value = '14:37'
$ws.Range("AB17").Value = '14:37'
$ws.Range("A18").Value = 130865719
$ws.Range("B18").Value = 79243
$ws.Range("Q18").Value = 447105
$ws.Range("R18").Value = 7043139
$ws.Range("S18").Value = 11
$ws.Range("Z18").Value = '15:35'
$ws.Range("AB18").Value = '15:35'
$ws.Range("AC18").Value = 'Med apothecier'
$ws.Range("B19").Value = 79243
$ws.Range("B20").Value = 79243
$ws.Range("B21").Value = 79243
$ws.Range("B22").Value = 91808
$ws.Range("B23").Value = 79243
$ws.Range("B24").Value = 79243
$ws.Range("A25").Value = 130865713
$ws.Range("B25").Value = 91804
$ws.Range("E25").Value = 1108
$ws.Range("F25").Value = 'Harticka'
$ws.Range("G25").Value = 'Pelloporus leporinus'
$ws.Range("H25").Value = '(Fr.) Krieglst.'
$ws.Range("Q25").Value = 447144
$ws.Range("R25").Value = 7043043
$ws.Range("S25").Value = 13
$ws.Range("Z25").Value = '15:18'
$ws.Range("AB25").Value = '15:18'
$ws.Range("A26").Value = 130865703
$ws.Range("B26").Value = 89193
$ws.Range("E26").Value = 510
$ws.Range("F26").Value = 'Doftskinn'
$ws.Range("G26").Value = 'Cystostereum murrayi'
$ws.Range("H26").Value = '(Berk. & M.A.Curtis.) Pouzar'
$ws.Range("Q26").Value = 447410
$ws.Range("R26").Value = 7042768
$ws.Range("S26").Value = 8
$ws.Range("Z26").Value = '14:20'
$ws.Range("AB26").Value = '14:20'
$ws.Range("A27").Value = 130865712
$ws.Range("B27").Value = 79243
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = 'Garnlav'
$ws.Range("G27").Value = 'Alectoria sarmentosa'
$ws.Range("H27").Value = '(Ach.) Ach.'
$ws.Range("Q27").Value = 447165
$ws.Range("R27").Value = 7043032
$ws.Range("S27").Value = 10
$ws.Range("Z27").Value = '15:13'
$ws.Range("AB27").Value = '15:13'
$ws.Range("B28").Value = 79243
$ws.Range("B29").Value = 57064
$ws.Range("B30").Value = 91804
$ws.Range("B31").Value = 79243

$ws.Range("AC16").ClearContents()
